$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (diff @@ -1223)
$ws.Range("H12").Value = 413.22223
$ws.Range("I12").Value = 339.875
$ws.Range("K12").Value = 339.875
$ws.Range("M12").Value = -169.875
# Row 28 (diff @@ -1989)
$ws.Range("H28").Value = 3068.697
$ws.Range("I28").Value = 3137.8845
$ws.Range("J28").Value = 2811.7144
$ws.Range("K28").Value = 3137.8845
$ws.Range("L28").Value = 2811.7144
$ws.Range("M28").Value = -2652.8845
$ws.Range("N28").Value = -3781.7144
# Row 53 (diff @@ -3214)
$ws.Range("H53").Value = 348.75
$ws.Range("I53").Value = 448.33334
$ws.Range("K53").Value = 448.33334
$ws.Range("M53").Value = 188.66666
# Row 63 (diff @@ -3713)
$ws.Range("H63").Value = 70135.5
$ws.Range("J63").Value = 70135.5
$ws.Range("L63").Value = 70135.5
$ws.Range("N63").Value = -71383.5
# Row 66 (diff @@ -3860)
$ws.Range("H66").Value = 70135.5
$ws.Range("J66").Value = 70135.5
$ws.Range("L66").Value = 210406.5
$ws.Range("N66").Value = -216646.5
# Row 88 (diff @@ -4971)
$ws.Range("H88").Value = 2031.7778
$ws.Range("I88").Value = 1099
$ws.Range("J88").Value = 2498.1667
$ws.Range("K88").Value = 1099
$ws.Range("L88").Value = 2498.1667
$ws.Range("M88").Value = -693
$ws.Range("N88").Value = -3310.1667
# Row 91 (diff @@ -5127)
$ws.Range("H91").Value = 2031.7778
$ws.Range("I91").Value = 1099
$ws.Range("J91").Value = 2498.1667
$ws.Range("K91").Value = 1099
$ws.Range("L91").Value = 2498.1667
$ws.Range("M91").Value = 305
$ws.Range("N91").Value = -5306.1667
# Row 101 (diff @@ -5629)
$ws.Range("H101").Value = 1343.2222
$ws.Range("I101").Value = 1156.4286
$ws.Range("J101").Value = 1997
$ws.Range("K101").Value = 3469.2858
$ws.Range("L101").Value = 5991
$ws.Range("M101").Value = -1847.2858
$ws.Range("N101").Value = -9235
# Row 130 (diff @@ -7053)
$ws.Range("H130").Value = 101000
$ws.Range("J130").Value = 101000
$ws.Range("L130").Value = 101000
$ws.Range("N130").Value = -111040
# Row 131 (diff @@ -7099)
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()
# Row 132 (diff @@ -7151)
$ws.Range("H132").Value = 501579.56
$ws.Range("J132").Value = 3336665.8
$ws.Range("L132").Value = 10009997.4
$ws.Range("N132").Value = -10015057.4
# Row 135 (diff @@ -7301)
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()
# Row 138 (diff @@ -7451)
$ws.Range("H138").Value = 2522.5925
$ws.Range("I138").Value = 862.2308
$ws.Range("J138").Value = 4064.3572
$ws.Range("K138").Value = 2586.6924
$ws.Range("L138").Value = 12193.0716
$ws.Range("M138").Value = 2553.3076
$ws.Range("N138").Value = -22473.0716
$ws = $wb.Worksheets.Item("BSM")
# Row 2 (diff @@ -14629)
$ws.Range("H2").Value = 60000
$ws.Range("J2").Value = 60000
$ws.Range("L2").Value = 60000
$ws.Range("N2").Value = -60226
# Row 69 (diff @@ -17873)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72 (diff @@ -18017)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 80 (diff @@ -18409)
$ws.Range("H80").Value = 638.1818
$ws.Range("I80").Value = 103.833336
$ws.Range("J80").Value = 1279.4
$ws.Range("K80").Value = 103.833336
$ws.Range("L80").Value = 1279.4
$ws.Range("M80").Value = 894.166664
$ws.Range("N80").Value = -3275.4
# Row 83 (diff @@ -18562)
$ws.Range("H83").Value = 638.1818
$ws.Range("I83").Value = 103.833336
$ws.Range("J83").Value = 1279.4
$ws.Range("K83").Value = 519.16668
$ws.Range("L83").Value = 6397
$ws.Range("M83").Value = 4472.83332
$ws.Range("N83").Value = -16381
# Row 99 (diff @@ -19367)
$ws.Range("H99").Value = 2018.2778
$ws.Range("I99").Value = 2004
$ws.Range("J99").Value = 2040.7142
$ws.Range("K99").Value = 2004
$ws.Range("L99").Value = 2040.7142
$ws.Range("M99").Value = -506
$ws.Range("N99").Value = -5036.7142
# Row 105 (diff @@ -19664)
$ws.Range("H105").Value = 4101.0625
$ws.Range("I105").Value = 4523.3335
$ws.Range("K105").Value = 4523.3335
$ws.Range("M105").Value = -2776.3335
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (diff @@ -22965)
$ws.Range("H31").Value = 7574.7095
$ws.Range("I31").Value = 3388.5
$ws.Range("J31").Value = 11022.177
$ws.Range("K31").Value = 3388.5
$ws.Range("L31").Value = 11022.177
$ws.Range("M31").Value = -3093.5
$ws.Range("N31").Value = -11612.177
# Row 34 (diff @@ -23115)
$ws.Range("H34").Value = 7574.7095
$ws.Range("I34").Value = 3388.5
$ws.Range("J34").Value = 11022.177
$ws.Range("K34").Value = 3388.5
$ws.Range("L34").Value = 11022.177
$ws.Range("M34").Value = -3186.5
$ws.Range("N34").Value = -11426.177
# Row 36 (diff @@ -23219)
$ws.Range("H36").Value = 2600
$ws.Range("I36").Value = 2600
$ws.Range("K36").Value = 2600
$ws.Range("M36").Value = -2212
# Row 40 (diff @@ -23418)
$ws.Range("H40").Value = 2600
$ws.Range("I40").Value = 2600
$ws.Range("K40").Value = 2600
$ws.Range("M40").Value = -2440
# Row 58 (diff @@ -24297)
$ws.Range("H58").Value = 2837.4614
$ws.Range("I58").Value = 2982.25
$ws.Range("K58").Value = 2982.25
$ws.Range("M58").Value = -2779.25
# Row 62 (diff @@ -24499)
$ws.Range("H62").Value = 12989.5
$ws.Range("I62").Value = 12486.875
$ws.Range("K62").Value = 12486.875
$ws.Range("M62").Value = -11862.875
# Row 65 (diff @@ -24649)
$ws.Range("H65").Value = 12989.5
$ws.Range("I65").Value = 12486.875
$ws.Range("K65").Value = 62434.375
$ws.Range("M65").Value = -59314.375
# Row 86 (diff @@ -25657)
$ws.Range("H86").Value = 3424.0908
$ws.Range("I86").Value = 3051.5
$ws.Range("J86").Value = 4417.6665
$ws.Range("K86").Value = 3051.5
$ws.Range("L86").Value = 4417.6665
$ws.Range("M86").Value = -1928.5
$ws.Range("N86").Value = -6663.6665
# Row 89 (diff @@ -25801)
$ws.Range("H89").Value = 3424.0908
$ws.Range("I89").Value = 3051.5
$ws.Range("J89").Value = 4417.6665
$ws.Range("K89").Value = 15257.5
$ws.Range("L89").Value = 22088.3325
$ws.Range("M89").Value = -9641.5
$ws.Range("N89").Value = -33320.3325
# Row 120 (diff @@ -27308)
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -47258
# Row 136 (diff @@ -28086)
$ws.Range("H136").Value = 2837.4614
$ws.Range("I136").Value = 2982.25
$ws.Range("K136").Value = 8946.75
$ws.Range("M136").Value = -6396.75
$ws = $wb.Worksheets.Item("CUL")
# Row 64 (diff @@ -31566)
$ws.Range("H64").Value = 1012
$ws.Range("I64").Value = 1012
$ws.Range("K64").Value = 3036
$ws.Range("M64").Value = -2766
# Row 67 (diff @@ -31707)
$ws.Range("H67").Value = 1012
$ws.Range("I67").Value = 1012
$ws.Range("K67").Value = 3036
$ws.Range("M67").Value = -2100
# Row 122 (diff @@ -34408)
$ws.Range("H122").Value = 1096.6666
$ws.Range("J122").Value = 1195
$ws.Range("L122").Value = 10755
$ws.Range("N122").Value = -15655
# Row 129 (diff @@ -34754)
$ws.Range("H129").Value = 8336.333000000001
$ws.Range("I129").Value = 2514.5
$ws.Range("K129").Value = 7543.5
$ws.Range("M129").Value = -2543.5
$ws = $wb.Worksheets.Item("GSM")
# Row 122 (diff @@ -41344)
$ws.Range("H122").Value = 1858.0834
$ws.Range("I122").Value = 1809.1
$ws.Range("J122").Value = 2103
$ws.Range("K122").Value = 5427.299999999999
$ws.Range("L122").Value = 6309
$ws.Range("M122").Value = -2977.299999999999
$ws.Range("N122").Value = -11209
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (diff @@ -42666)
$ws.Range("H7").Value = 3999.625
$ws.Range("I7").Value = 4073
$ws.Range("K7").Value = 4073
$ws.Range("M7").Value = -3961
# Row 46 (diff @@ -44598)
$ws.Range("H46").Value = 1169
$ws.Range("I46").Value = 1126.4546
$ws.Range("K46").Value = 1126.4546
$ws.Range("M46").Value = -938.4546
# Row 100 (diff @@ -47196)
$ws.Range("H100").Value = 2825.8096
$ws.Range("I100").Value = 2686.625
$ws.Range("J100").Value = 2911.4614
$ws.Range("K100").Value = 2686.625
$ws.Range("L100").Value = 2911.4614
$ws.Range("M100").Value = -2145.625
$ws.Range("N100").Value = -3993.4614
# Row 126 (diff @@ -48464)
$ws.Range("H126").Value = 3999.625
$ws.Range("I126").Value = 4073
$ws.Range("K126").Value = 12219
$ws.Range("M126").Value = -9749
# Row 132 (diff @@ -48758)
$ws.Range("H132").Value = 3254.077
$ws.Range("I132").Value = 2786.2856
$ws.Range("K132").Value = 8358.856800000001
$ws.Range("M132").Value = -5828.856800000001
$ws = $wb.Worksheets.Item("WVR")
# Row 41 (diff @@ -51256)
$ws.Range("H41").Value = 14332
$ws.Range("J41").Value = 14498
$ws.Range("L41").Value = 14498
$ws.Range("N41").Value = -15278
# Row 107 (diff @@ -54487)
$ws.Range("H107").Value = 2063
$ws.Range("I107").Value = 2166.1667
$ws.Range("K107").Value = 6498.500100000001
$ws.Range("M107").Value = -4578.500100000001
# Row 126 (diff @@ -55409)
$ws.Range("H126").Value = 3566.0833
$ws.Range("I126").Value = 3136.4736
$ws.Range("K126").Value = 9409.4208
$ws.Range("M126").Value = -6939.4208
# Row 132 (diff @@ -55697)
$ws.Range("H132").Value = 5880.364
$ws.Range("I132").Value = 5411.0625
$ws.Range("K132").Value = 16233.1875
$ws.Range("M132").Value = -13703.1875
